$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.331.41'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '2.527.55'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''323.70'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '''109.57'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.558'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +4.25%  '
$ws.Range("D10").Value = '''40.71'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +4.09%  '
$ws.Range("D11").Value = '''20.37'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +10.97%  '
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("E14").Value = '  +1.51%  '
$ws.Range("D15").Value = '2.924.64'
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").Value = '2.533.55'
$ws.Range("E16").Value = '  +1.17%  '
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").Value = '48.180.63'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").Value = '''13.32'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("D20").Value = '''6.64'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").Value = '''72.54'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("D24").Value = '''269.75'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +8.84%  '
$ws.Range("D25").Value = '''2.58'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("D26").Value = '''26.29'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '''10.19'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("E29").Value = '  +5.85%  '
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("E31").Value = '  -8.48%  '
$ws.Range("D32").Value = '''49.79'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").Value = '''20.07'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").Value = '''22.40'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").Value = '''118.60'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -1.95%  '
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("D45").Value = '2.015.19'
$ws.Range("E45").Value = '  +1.20%  '
$ws.Range("E46").Value = '  +3.47%  '
$ws.Range("D47").Value = '''1.89'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +6.13%  '
$ws.Range("E48").Value = '  -1.24%  '
$ws.Range("D49").Value = '''9.15'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '''80.03'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +2.55%  '
